$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("VT-SaleVoid-DualCF-Generic")
$ws.Range("B2").Value = "Thu Nov 20 20:57:37 IST 2025"
$ws.Range("B3").Value = "Thu Nov 20 20:58:40 IST 2025"
$ws.Range("B4").Value = "Thu Nov 20 20:59:48 IST 2025"
$ws.Range("B5").Value = "Thu Nov 20 21:00:50 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleVoid-NoCF-Generic")
$ws.Range("B2").Value = "Thu Nov 20 16:43:32 IST 2025"
$ws.Range("D2").Value = "Thu Nov 20 21:01:53 IST 2025"
$ws.Range("B3").Value = "Thu Nov 20 16:44:31 IST 2025"
$ws.Range("D3").Value = "Thu Nov 20 21:02:48 IST 2025"
$ws.Range("B4").Value = "Thu Nov 20 16:45:32 IST 2025"
$ws.Range("D4").Value = "Thu Nov 20 21:03:42 IST 2025"
$ws.Range("B5").Value = "Thu Nov 20 16:46:35 IST 2025"
$ws.Range("D5").Value = "Thu Nov 20 21:04:38 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleVoid-SingleCF-Generic")
$ws.Range("B2").Value = "Thu Nov 20 21:05:41 IST 2025"
$ws.Range("B3").Value = "Thu Nov 20 21:06:38 IST 2025"
$ws.Range("B4").Value = "Thu Nov 20 21:07:31 IST 2025"
$ws.Range("B5").Value = "Thu Nov 20 21:08:30 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleCredit-DualCF-Generic")
$ws.Range("B2").Value = "Thu Nov 20 20:47:03 IST 2025"
$ws.Range("B3").Value = "Thu Nov 20 20:48:00 IST 2025"
$ws.Range("B4").Value = "Thu Nov 20 20:49:00 IST 2025"
$ws.Range("B5").Value = "Thu Nov 20 20:49:49 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleCredit-NoCF-Generic")
$ws.Range("B2").Value = "Thu Nov 20 20:50:44 IST 2025"
$ws.Range("B3").Value = "Thu Nov 20 20:51:28 IST 2025"
$ws.Range("B4").Value = "Thu Nov 20 20:52:21 IST 2025"
$ws.Range("B5").Value = "Thu Nov 20 20:53:11 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleCredit-SingleCF-Generic")
$ws.Range("B2").Value = "Thu Nov 20 20:54:12 IST 2025"
$ws.Range("C2").Value = "Fail"
$ws.Range("B3").Value = "Thu Nov 20 20:55:07 IST 2025"
$ws.Range("C3").Value = "Fail"
$ws.Range("B4").Value = "Thu Nov 20 20:56:01 IST 2025"
$ws.Range("C4").Value = "Fail"
$ws.Range("B5").Value = "Thu Nov 20 20:56:44 IST 2025"
$ws.Range("C5").Value = "Fail"

$ws = $wb.Worksheets.Item("VT-AuthCapCredit-Generic")
$ws.Range("B2").Value = "Thu Nov 20 15:36:44 IST 2025"
$ws.Range("D2").Value = "Thu Nov 20 22:16:12 IST 2025"
$ws.Range("B3").Value = "Thu Nov 20 15:38:06 IST 2025"
$ws.Range("D3").Value = "Thu Nov 20 22:17:33 IST 2025"
$ws.Range("B4").Value = "Thu Nov 20 15:39:29 IST 2025"
$ws.Range("C4").Value = "Fail"
$ws.Range("D4").Value = "Thu Nov 20 22:18:51 IST 2025"
$ws.Range("B5").Value = "Thu Nov 20 15:40:48 IST 2025"
$ws.Range("C5").Value = "Fail"
$ws.Range("D5").Value = "Thu Nov 20 22:20:01 IST 2025"
$ws.Range("B6").Value = "Thu Nov 20 15:42:07 IST 2025"
$ws.Range("D6").Value = "Thu Nov 20 22:21:29 IST 2025"
$ws.Range("B7").Value = "Thu Nov 20 15:43:37 IST 2025"
$ws.Range("D7").Value = "Thu Nov 20 22:22:55 IST 2025"

$ws = $wb.Worksheets.Item("VT-AuthCapVoid-Generic")
$ws.Range("B2").Value = "Thu Nov 20 21:12:50 IST 2025"
$ws.Range("D2").Value = "Thu Nov 20 19:46:16 IST 2025"
$ws.Range("B3").Value = "Thu Nov 20 21:14:29 IST 2025"
$ws.Range("D3").Value = "Thu Nov 20 19:47:47 IST 2025"
$ws.Range("B4").Value = "Thu Nov 20 21:16:08 IST 2025"
$ws.Range("D4").Value = "Thu Nov 20 19:49:01 IST 2025"
$ws.Range("B5").Value = "Thu Nov 20 21:17:35 IST 2025"
$ws.Range("D5").Value = "Thu Nov 20 19:50:28 IST 2025"
$ws.Range("B6").Value = "Thu Nov 20 21:18:59 IST 2025"
$ws.Range("D6").Value = "Thu Nov 20 19:51:54 IST 2025"
$ws.Range("B7").Value = "Thu Nov 20 21:20:25 IST 2025"
$ws.Range("D7").Value = "Thu Nov 20 19:53:16 IST 2025"

$ws = $wb.Worksheets.Item("VT-ManualAuthCapture-Generic")
$ws.Range("B2").Value = "Thu Nov 20 20:40:33 IST 2025"
$ws.Range("B3").Value = "Thu Nov 20 20:41:30 IST 2025"
$ws.Range("B4").Value = "Thu Nov 20 20:42:24 IST 2025"
$ws.Range("B5").Value = "Thu Nov 20 20:43:54 IST 2025"
$ws.Range("B6").Value = "Thu Nov 20 20:45:02 IST 2025"
$ws.Range("B7").Value = "Thu Nov 20 20:46:10 IST 2025"
